# Generate Report for Handback
# Rename the two e2e test markdown files tracked by this handback-status
# report and refresh the handback timestamps / correspond-handback xliff
# file names that the CI run regenerated.

$wb = $excel.ActiveWorkbook

# ---- new identifiers -------------------------------------------------
$newName1 = "c825b56e-174d-4a96-a1e6-bec693c5b85b.md"
$newName2 = "ffffdad39318-0ab9-4f71-81e2-dc7c6c1c5b4b.md"
$newPath1 = "e2e\c825b56e-174d-4a96-a1e6-bec693c5b85b.md"
$newPath2 = "e2e\ffffdad39318-0ab9-4f71-81e2-dc7c6c1c5b4b.md"

$overviewDate = "2016-08-12 05:05:51"

$zhcnXlf    = "c825b56e-174d-4a96-a1e6-bec693c5b85b.40e450010f2f9e154ec5e3d990e094f46020c88d.zh-cn.xlf"
$zhcnHDate  = "2016-08-12 05:05:46"
$zhcnKDate  = "2016-08-12 05:06:14"

$dedeXlf    = "c825b56e-174d-4a96-a1e6-bec693c5b85b.40e450010f2f9e154ec5e3d990e094f46020c88d.de-de.xlf"
$dedeKDate  = "2016-08-12 05:06:21"

# ---- Overview sheet ----------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newName1
$wsOverview.Range("A3").Value = $newName2
$wsOverview.Range("G2").Value = $overviewDate
$wsOverview.Range("G3").Value = $overviewDate

$wsOverview.Hyperlinks.Item(1).TextToDisplay = $newPath1
$wsOverview.Hyperlinks.Item(2).TextToDisplay = $newPath2

# ---- zh-cn sheet ---------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newName1
$wsZhCn.Range("I2").Value = $newName1
$wsZhCn.Range("G2").Value = $zhcnXlf
$wsZhCn.Range("J2").Value = $zhcnXlf
$wsZhCn.Range("H2").Value = $zhcnHDate
$wsZhCn.Range("K2").Value = $zhcnKDate

$wsZhCn.Range("A3").Value = $newName2
$wsZhCn.Range("I3").Value = $newName2
$wsZhCn.Range("G3").Value = $zhcnXlf
$wsZhCn.Range("J3").Value = $zhcnXlf
$wsZhCn.Range("H3").Value = $zhcnHDate
$wsZhCn.Range("K3").Value = $zhcnKDate

$wsZhCn.Hyperlinks.Item(1).TextToDisplay = $newName1
$wsZhCn.Hyperlinks.Item(2).TextToDisplay = $newName1
$wsZhCn.Hyperlinks.Item(3).TextToDisplay = $newName2
$wsZhCn.Hyperlinks.Item(4).TextToDisplay = $newName2

# ---- de-de sheet ---------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newName1
$wsDeDe.Range("I2").Value = $newName1
$wsDeDe.Range("G2").Value = $dedeXlf
$wsDeDe.Range("J2").Value = $dedeXlf
$wsDeDe.Range("H2").Value = $overviewDate
$wsDeDe.Range("K2").Value = $dedeKDate

$wsDeDe.Range("A3").Value = $newName2
$wsDeDe.Range("I3").Value = $newName2
$wsDeDe.Range("G3").Value = $dedeXlf
$wsDeDe.Range("J3").Value = $dedeXlf
$wsDeDe.Range("H3").Value = $overviewDate
$wsDeDe.Range("K3").Value = $dedeKDate

$wsDeDe.Hyperlinks.Item(1).TextToDisplay = $newName1
$wsDeDe.Hyperlinks.Item(2).TextToDisplay = $newName1
$wsDeDe.Hyperlinks.Item(3).TextToDisplay = $newName2
$wsDeDe.Hyperlinks.Item(4).TextToDisplay = $newName2
